$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (Changed) date, stored as an Excel serial
# date number. Rows 2 through 176 currently contain 45204 and need to be
# updated to 45207.
for ($row = 2; $row -le 176; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45207
    }
}
